$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 4477
$ws.Range("C9").Value = 4226
$ws.Range("D9").Value = 4451
$ws.Range("E9").Value = 4477
$ws.Range("F9").Value = 4477
